$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A handful of blank (but bold-formatted) rows were left behind further down
# the sheet by the query helper's day-by-day paste runs. Stake out their
# style first so it claims the earlier style slot (matches how the workbook
# that produced this diff ended up laying out its cellXfs table).
$boldMaster = $ws.Cells.Item(64, 1)
$boldMaster.Font.Bold = $true
$boldMaster.Font.Name = "Lucida Grande"
$boldMaster.Font.Size = 11
$boldMaster.Font.Color = 0

# New language/code pairs pulled in from a few days of additional data
$newRows = @(
    @("cs","Czech"),
    @("da","Danish"),
    @("de","German"),
    @("de-DE","German (Germany)"),
    @("en-CA","English (Canada)"),
    @("en-GB","English (United Kingdom)"),
    @("en-US","English (United States)"),
    @("es-ES","Spanish (Spain)"),
    @("es-MX","Spanish (Mexico)"),
    @("fa-IR","Persian (Iran)"),
    @("fi","Finnish"),
    @("fr-FR","French (France)"),
    @("hi","Hindi"),
    @("ms","Malay"),
    @("nl","Dutch"),
    @("nl-NL","Dutch (Netherands)"),
    @("no","Norwegian"),
    @("pt-PT","Portugese (Portugal)"),
    @("sv","Swedish"),
    @("yue-HK","Cantonese (Hong Kong)"),
    @("zh-CN","Chinese (China)"),
    @("zh-Hans","Chinese (Simplified)"),
    @("zh-Hant","Chinese (Traditional)"),
    @("zh-HK","Chinese (Hong Kong)"),
    @("zh-SG","Chinese (Singapore)")
)

$startRow = 22
$r = $startRow
foreach ($pair in $newRows) {
    $ws.Cells.Item($r, 1).Value = $pair[0]
    $ws.Cells.Item($r, 2).Value = $pair[1]
    $r++
}
$endRow = $r - 1

# the pasted-in codes came through with their source formatting (11pt Lucida
# Grande, black) on the code column only -- reproduce that with one "master"
# cell then fan it out with copy/paste so we don't touch the language column
$fontMaster = $ws.Cells.Item($startRow, 1)
$fontMaster.Font.Name = "Lucida Grande"
$fontMaster.Font.Size = 11
$fontMaster.Font.Color = 0
$fontMaster.Copy()
$ws.Range($ws.Cells.Item($startRow + 1, 1), $ws.Cells.Item($endRow, 1)).PasteSpecial(-4122)

# resort the table by language code, keeping the header in place
$fullRange = $ws.Range("A1:B" + $endRow)
$fullRange.Sort($ws.Range("A1"), 1, $null, $null, 1, $null, $null, 1)

$lastDataRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

$dataRange = $ws.Range("A1:B" + $lastDataRow)
$dataRange.AutoFilter()

$filterName = $ws.Names.Add("_xlnm._FilterDatabase", "=" + $ws.Name + "!`$A`$1:`$B`$" + $lastDataRow)
$filterName.Visible = $false

# fan the bold-blank style out to the rest of the leftover block
$boldMaster.Copy()
$ws.Range($ws.Cells.Item(65, 1), $ws.Cells.Item(67, 1)).PasteSpecial(-4122)
